# Updated main GSC export data
#
# The "Chart" sheet's daily export table drops its oldest day (2025-10-31,
# previously row 2) -- every remaining row shifts up by one -- and gains
# three new trailing days (2026-01-25, 2026-01-26, 2026-01-27).
#
# The "Critical issues" sheet's Pages counts are refreshed for four reasons.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Chart" sheet: drop the 2025-10-31 row (row 2) and shift everything
# else up by one. Using EntireRow delete (rather than rewriting every
# cell) keeps the existing text/number cell types and default (General)
# style intact for all the untouched rows.
# ---------------------------------------------------------------------
$wsChart = $wb.Worksheets.Item("Chart")
$wsChart.Rows.Item(2).Delete()

# Append the three new trailing days. Column A holds a literal date-like
# text string (not a real date value) in this workbook, so we briefly
# force Text format before assigning the value (otherwise Excel parses
# "2026-01-25" as a date serial), then clear the formatting again so the
# cell ends up back on the default/General style -- matching every other
# cell in the sheet -- while keeping the text content.
function Set-DateText($ws, $cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-DateText $wsChart "A87" "2026-01-25"
$wsChart.Range("B87").Value = 304
$wsChart.Range("C87").Value = 223
$wsChart.Range("D87").Value = 70

Set-DateText $wsChart "A88" "2026-01-26"
$wsChart.Range("B88").Value = 304
$wsChart.Range("C88").Value = 223
$wsChart.Range("D88").Value = 85

Set-DateText $wsChart "A89" "2026-01-27"
$wsChart.Range("B89").Value = 304
$wsChart.Range("C89").Value = 223
$wsChart.Range("D89").Value = 91

# ---------------------------------------------------------------------
# "Critical issues" sheet: refresh the Pages counts for four reasons.
# ---------------------------------------------------------------------
$wsCritical = $wb.Worksheets.Item("Critical issues")
$wsCritical.Range("D2").Value = 96
$wsCritical.Range("D4").Value = 60
$wsCritical.Range("D5").Value = 26
$wsCritical.Range("D10").Value = 19
